$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 10; $r++) {
    foreach ($col in @("A", "B")) {
        $cell = $ws.Range("$col$r")
        $val = $cell.Value2
        if ($val -ne $null) {
            $cell.Value2 = $val.Replace(":", "=")
        }
    }
}
